# New crime data collected - weekly CompStat (cs-en-us-pbsi) refresh.
# Updates the report volume/number, the reporting week dates, and the
# weekly/28-day/YTD/2-year crime-complaint figures (rows 14-31) to the
# newly collected numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Masthead: bump the Volume/Number and the "Report Covering the Week"
#    date range by one week.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/5/2024  Through  8/11/2024"

# ---------------------------------------------------------------------
# 2) Weekly crime-complaint grid (columns C:N, rows 14-31).
#    Each hashtable below lists only the columns that actually change
#    for that row; values are written straight into the existing cell
#    formatting (counts, % changes, etc. keep their current style).
# ---------------------------------------------------------------------
$rowUpdates = @(
    @{ Row = 14; Cells = @{ G = 3;  H = -33.333333333333 } }
    @{ Row = 15; Cells = @{ D = 1;  E = 0;                  F = 6;   G = 2;
                            H = 200; I = 38; J = 24;
                            K = 58.333333333333; L = 111.111111111111;
                            M = 18.75;            N = -29.629629629629 } }
    @{ Row = 16; Cells = @{ C = 5;  D = 9;  E = -44.444444444444;
                            F = 21; G = 29; H = -27.586206896551;
                            I = 190; J = 199; K = -4.522613065326;
                            L = 27.516778523489; M = -22.764227642276;
                            N = -76.220275344180 } }
    @{ Row = 17; Cells = @{ C = 18; D = 15; E = 20;
                            F = 75; G = 73; H = 2.739726027397;
                            I = 590; J = 597; K = -1.172529313232;
                            L = 18; M = 94.078947368421;
                            N = -19.618528610354 } }
    @{ Row = 18; Cells = @{ C = 8;  D = 9;  E = -11.111111111111;
                            F = 31; G = 37; H = -16.216216216216;
                            I = 190; J = 234; K = -18.803418803418;
                            L = 12.426035502958; M = -46.778711484593;
                            N = -90.969581749049 } }
    @{ Row = 19; Cells = @{ C = 25; D = 32; E = -21.875;
                            F = 86; G = 147; H = -41.496598639455;
                            I = 934; J = 973; K = -4.008221993833;
                            L = 8.857808857808; M = 58.573853989813;
                            N = -6.506506506506 } }
    @{ Row = 20; Cells = @{ D = 15; E = -46.666666666666;
                            F = 32; G = 56; H = -42.857142857142;
                            I = 216; J = 281; K = -23.131672597864;
                            L = -16.279069767441; M = -0.917431192660;
                            N = -92.727272727272 } }
    @{ Row = 21; Cells = @{ C = 65; D = 81; E = -19.753086419753;
                            F = 253; G = 347; H = -27.089337175792;
                            I = 2162; J = 2323; K = -6.930693069306;
                            L = 10.418794688457; M = 23.190883190883;
                            N = -71.834288692027 } }
    @{ Row = 23; Cells = @{ C = 4;  D = 4;  E = 0;
                            F = 14; H = 27.272727272727;
                            I = 72; J = 95; K = -24.210526315789;
                            L = 10.769230769230; M = 89.473684210526 } }
    @{ Row = 24; Cells = @{ C = 51; D = 83; E = -38.554216867469;
                            F = 257; G = 361; H = -28.808864265928;
                            I = 2508; J = 2593; K = -3.278056305437;
                            L = 10.728476821192; M = 9.711286089238 } }
    @{ Row = 25; Cells = @{ C = 29; D = 36; E = -19.444444444444;
                            F = 118; G = 157; H = -24.840764331210;
                            I = 1327; J = 1216; K = 9.128289473684;
                            L = 62.026862026862 } }
    @{ Row = 26; Cells = @{ C = 36; D = 46; E = -21.739130434782;
                            F = 143; G = 152; H = -5.921052631578;
                            I = 1164; J = 1102; K = 5.626134301270;
                            L = 13.560975609756; M = -5.596107055961 } }
    @{ Row = 27; Cells = @{ C = 3;  D = 2;  E = 50;
                            F = 8;  H = 166.666666666667;
                            I = 65; J = 41; K = 58.536585365853;
                            L = 62.5 } }
    @{ Row = 28; Cells = @{ C = 5;  E = 150; G = 10;
                            H = 130; I = 130; J = 127;
                            K = 2.362204724409; L = 21.495327102803 } }
    @{ Row = 29; Cells = @{ D = 1;  E = -100; G = 3;
                            H = 0;  J = 23;
                            K = -52.173913043478; L = -57.692307692307;
                            N = -84.057971014492 } }
    @{ Row = 30; Cells = @{ D = 1;  E = -100; G = 3;
                            H = -33.333333333333; J = 21;
                            K = -52.380952380952; L = -52.380952380952;
                            N = -83.050847457627 } }
    @{ Row = 31; Cells = @{ L = -17.647058823529 } }
)

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $ws.Range("$col$row").Value = $update.Cells[$col]
    }
}

# ---------------------------------------------------------------------
# 3) A handful of cells flip between a numeric value and the sheet's
#    "no data" text placeholders ("0" / "***.*"). Those placeholders
#    must stay genuine text (like the neighbouring label cells) rather
#    than become a number, and a couple of previously-placeholder cells
#    now need to go back to being real numbers, so their formatting is
#    fixed up explicitly here.
# ---------------------------------------------------------------------

# Cells that must display as a real, calculated number again (restore
# the normal count/percent formatting from a same-column donor cell).
$numericFixups = @(
    @{ Cell = "D15"; Value = 1;  FormatFrom = "C16" }
    @{ Cell = "E15"; Value = 0;  FormatFrom = "H16" }
    @{ Cell = "D27"; Value = 2;  FormatFrom = "C16" }
    @{ Cell = "E27"; Value = 50; FormatFrom = "H16" }
)
foreach ($fix in $numericFixups) {
    $dst = $ws.Range($fix.Cell)
    $dst.Value = $fix.Value
    $ws.Range($fix.FormatFrom).Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# Cells that must become the literal placeholder text "0" / "***.*"
# (kept as text, formatted like the row's label cell in column A).
$textPlaceholders = @(
    @{ Cell = "C14"; Value = "0" }
    @{ Cell = "D14"; Value = "0" }
    @{ Cell = "E14"; Value = "***.*" }
    @{ Cell = "C29"; Value = "0" }
    @{ Cell = "C30"; Value = "0" }
    @{ Cell = "D31"; Value = "0" }
    @{ Cell = "E31"; Value = "***.*" }
)
foreach ($ph in $textPlaceholders) {
    $dst = $ws.Range($ph.Cell)
    $row = $dst.Row
    $dst.Value = "'" + $ph.Value
    $ws.Range("A$row").Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = $false
}

Write-Host "CompStat weekly data refreshed."
